$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I21").Value = 20
$ws.Range("H23").Value = 30
$ws.Range("I23").Value = 250
$ws.Range("D24").Value = 1
$ws.Range("H24").Value = 40
$ws.Range("I24").Value = 20

$excel.Calculate()
